$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Select()

# Capture the width of the column immediately to the left (M) so the
# newly inserted column (N) inherits it, exactly like Excel does when you
# insert a column via the UI.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"); this shifts the
# former N/O/P ("Late", heading/"Outstanding" label, "Outstanding")
# columns one place to the right, to O/P/Q.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftWidth

# Leave the selection where the user ended up after the edit.
$ws.Range("S8").Select()
